# The workbook tracks daily price observations for "Poroto verde" at the
# "Vega Modelo de Temuco". A new observation (dated 2022-06-03) needs to be
# inserted right before the existing row 39 (dated 2020-11-26), shifting the
# old row 39 and everything below it down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 39; this pushes old rows 39..139 down to 40..140.
$ws.Rows.Item(39).Insert()

# The new row 39 repeats the (unchanged) descriptive columns from the row
# that used to sit there (now row 40) and supplies the new observation's
# Fecha / Volumen / Precio / Unidad values.
$ws.Range("A39").Value = 10
$ws.Range("B39").Value = "Vega Modelo de Temuco"
$ws.Range("C39").Value = "La Araucanía"
$ws.Range("D39").Value = 44715
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = 100112031
$ws.Range("G39").Value = "Poroto verde"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 55
$ws.Range("K39").Value = 25000
$ws.Range("L39").Value = 25000
$ws.Range("M39").Value = 25000
$ws.Range("N39").Value = "$/malla 25 kilos"
$ws.Range("O39").Value = "Región del Maule"
$ws.Range("P39").Value = 1000
$ws.Range("Q39").Value = 25
$ws.Range("R39").Value = "Hortaliza"
